$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 7085873
$ws.Range("I8").Value = 9018194
$ws.Range("J8").Value = 696
$ws.Range("K8").Value = 27054582
$ws.Range("L8").Value = 2088
$ws.Range("M8").Value = -27054443
$ws.Range("N8").Value = -2366
$ws.Range("H17").Value = 528009
$ws.Range("J17").Value = 528009
$ws.Range("L17").Value = 1584027
$ws.Range("N17").Value = -1584363
$ws.Range("H39").Value = 672.61536
$ws.Range("I39").Value = 60.166668
$ws.Range("J39").Value = 1197.5714
$ws.Range("K39").Value = 180.500004
$ws.Range("L39").Value = 3592.7142
$ws.Range("M39").Value = 115.499996
$ws.Range("N39").Value = -4184.7142
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()
$ws.Range("H97").Value = 3253.75
$ws.Range("J97").Value = 3253.75
$ws.Range("L97").Value = 9761.25
$ws.Range("N97").Value = -10753.25
$ws.Range("H106").Value = 8355.388999999999
$ws.Range("I106").Value = 2105.5454
$ws.Range("K106").Value = 2105.5454
$ws.Range("M106").Value = -1474.5454
$ws.Range("H112").Value = 1536.9333
$ws.Range("H116").Value = 4247.909
$ws.Range("I116").Value = 3992.375
$ws.Range("K116").Value = 3992.375
$ws.Range("M116").Value = -550.375
$ws.Range("H138").Value = 2577.2307
$ws.Range("I138").Value = 1986.5
$ws.Range("J138").Value = 2710.9812
$ws.Range("K138").Value = 5959.5
$ws.Range("L138").Value = 8132.943600000001
$ws.Range("M138").Value = -819.5
$ws.Range("N138").Value = -18412.9436
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 3970.4443
$ws.Range("I21").Value = 3305.5
$ws.Range("J21").Value = 4502.4
$ws.Range("K21").Value = 3305.5
$ws.Range("L21").Value = 4502.4
$ws.Range("M21").Value = -2931.5
$ws.Range("N21").Value = -5250.4
$ws.Range("H32").Value = 4527.909
$ws.Range("I32").Value = 3534.507
$ws.Range("K32").Value = 3534.507
$ws.Range("M32").Value = -3247.507
$ws.Range("H61").Value = 1807.6923
$ws.Range("I61").Value = 1715.5834
$ws.Range("K61").Value = 1715.5834
$ws.Range("M61").Value = -1503.5834
$ws.Range("H136").Value = 1807.6923
$ws.Range("I136").Value = 1715.5834
$ws.Range("K136").Value = 5146.7502
$ws.Range("M136").Value = -2596.7502
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1839.3846
$ws.Range("I86").Value = 1960.5
$ws.Range("K86").Value = 1960.5
$ws.Range("M86").Value = -837.5
$ws.Range("H89").Value = 1839.3846
$ws.Range("I89").Value = 1960.5
$ws.Range("K89").Value = 9802.5
$ws.Range("M89").Value = -4186.5
$ws.Range("H134").Value = 89506.06
$ws.Range("I134").Value = 97977.23
$ws.Range("J134").Value = 1970.6666
$ws.Range("K134").Value = 293931.69
$ws.Range("L134").Value = 5911.9998
$ws.Range("M134").Value = -291396.69
$ws.Range("N134").Value = -10981.9998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 152.14285
$ws.Range("I17").Value = 57
$ws.Range("J17").Value = 390
$ws.Range("K17").Value = 171
$ws.Range("L17").Value = 1170
$ws.Range("M17").Value = -2
$ws.Range("N17").Value = -1508
$ws.Range("H23").Value = 235.57143
$ws.Range("I23").Value = 267.8
$ws.Range("J23").Value = 155
$ws.Range("K23").Value = 803.4000000000001
$ws.Range("L23").Value = 465
$ws.Range("M23").Value = -568.4000000000001
$ws.Range("N23").Value = -935
$ws.Range("H113").Value = 1810.72
$ws.Range("I113").Value = 501.5
$ws.Range("K113").Value = 1504.5
$ws.Range("M113").Value = 665.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2500548.5
$ws.Range("I3").Value = 2000077.6
$ws.Range("J3").Value = 3334666.8
$ws.Range("K3").Value = 2000077.6
$ws.Range("L3").Value = 3334666.8
$ws.Range("M3").Value = -1999961.6
$ws.Range("N3").Value = -3334898.8
$ws.Range("I107").Value = 660.6
$ws.Range("J107").Value = 62506820
$ws.Range("K107").Value = 660.6
$ws.Range("L107").Value = 62506820
$ws.Range("M107").Value = 1259.4
$ws.Range("N107").Value = -62510660
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H113").Value = 2933.9167
$ws.Range("I113").Value = 2022.2222
$ws.Range("K113").Value = 2022.2222
$ws.Range("M113").Value = 147.7778000000001
$ws.Range("H122").Value = 3829.5
$ws.Range("I122").Value = 3575.5
$ws.Range("J122").Value = 3998.8333
$ws.Range("K122").Value = 10726.5
$ws.Range("L122").Value = 11996.4999
$ws.Range("M122").Value = -8276.5
$ws.Range("N122").Value = -16896.4999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 28661.06
$ws.Range("I7").Value = 21864.857
$ws.Range("J7").Value = 33668.79
$ws.Range("K7").Value = 21864.857
$ws.Range("L7").Value = 33668.79
$ws.Range("M7").Value = -21752.857
$ws.Range("N7").Value = -33892.79
$ws.Range("H16").Value = 704.3333
$ws.Range("I16").Value = 750.63635
$ws.Range("J16").Value = 195
$ws.Range("K16").Value = 750.63635
$ws.Range("L16").Value = 195
$ws.Range("M16").Value = -580.63635
$ws.Range("N16").Value = -535
$ws.Range("H122").Value = 78112.19
$ws.Range("I122").Value = 108677.52
$ws.Range("J122").Value = 5519.5
$ws.Range("K122").Value = 326032.56
$ws.Range("L122").Value = 16558.5
$ws.Range("M122").Value = -323582.56
$ws.Range("N122").Value = -21458.5
$ws.Range("H126").Value = 28661.06
$ws.Range("I126").Value = 21864.857
$ws.Range("J126").Value = 33668.79
$ws.Range("K126").Value = 65594.571
$ws.Range("L126").Value = 101006.37
$ws.Range("M126").Value = -63124.571
$ws.Range("N126").Value = -105946.37
$ws.Range("H132").Value = 4445.657
$ws.Range("I132").Value = 3848.889
$ws.Range("J132").Value = 6459.75
$ws.Range("K132").Value = 11546.667
$ws.Range("L132").Value = 19379.25
$ws.Range("M132").Value = -9016.667000000001
$ws.Range("N132").Value = -24439.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3024
$ws.Range("I81").Value = 3024
$ws.Range("K81").Value = 6048
$ws.Range("M81").Value = -4987
$ws.Range("H84").Value = 3024
$ws.Range("I84").Value = 3024
$ws.Range("K84").Value = 30240
$ws.Range("M84").Value = -24936
$ws.Range("H98").Value = 30590
$ws.Range("J98").Value = 30590
$ws.Range("L98").Value = 30590
$ws.Range("N98").Value = -36580
$ws.Range("H113").Value = 415.64706
$ws.Range("I113").Value = 391.27274
$ws.Range("K113").Value = 1173.81822
$ws.Range("M113").Value = 996.1817799999999
$ws.Range("H122").Value = 1214.7102
$ws.Range("I122").Value = 1150.4182
$ws.Range("J122").Value = 1467.2858
$ws.Range("K122").Value = 3451.2546
$ws.Range("L122").Value = 4401.857400000001
$ws.Range("M122").Value = -1001.2546
$ws.Range("N122").Value = -9301.857400000001
